# Commit: "Add files via upload"
#
# Net effect of the OOXML diff: slide 17 ("Slide 17", an empty filler
# slide with no notes content) is removed from the deck. Every part
# after it (slides/notesSlides 18-25, their relationship ids, the
# section's p14:sldId list, and the cached "slide number" text fields
# on the following notes pages) merely shifts down by one position as
# an automatic side effect of PowerPoint re-saving the package after
# the deletion - none of that reflects an independent content edit, so
# the single COM action below reproduces the authored change.

$p = $ppt.ActivePresentation

# "Slide 17" is the 17th slide in presentation order (rId18 / sldId 272).
$p.Slides.Item(17).Delete()
